$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.551034
$ws.Range("H2").Value = 1.653102
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.2433346666666667
$ws.Range("N2").Value = 0.730004
$ws.Range("O2").Value = 0.0006441701120846545
$ws.Range("P2").Value = 0.0006447102349388058
$ws.Range("Q2").Value = 0.134085674712
$ws.Range("R2").Value = 1.206771072408
$ws.Range("S2").Value = 0.0006441701120846545
$ws.Range("T2").Value = 0.0006447102349388058

# Row 3
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.551034
$ws.Range("H3").Value = 1.653102
$ws.Range("O3").Value = 0.000824430864760877
$ws.Range("P3").Value = 0.0008251221323986817
$ws.Range("Q3").Value = 0.171607416552
$ws.Range("R3").Value = 1.544466748968
$ws.Range("S3").Value = 0.000824430864760877
$ws.Range("T3").Value = 0.0008251221323986817

# Row 4
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.551034
$ws.Range("H4").Value = 1.653102
$ws.Range("M4").Value = 318.203888
$ws.Range("N4").Value = 954.611664
$ws.Range("O4").Value = 0.8423684015377977
$ws.Range("P4").Value = 0.843074709416338
$ws.Range("Q4").Value = 175.341161220192
$ws.Range("R4").Value = 1578.070450981728
$ws.Range("S4").Value = 0.8423684015377977
$ws.Range("T4").Value = 0.843074709416338

# Row 5
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.551034
$ws.Range("H5").Value = 1.653102
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.5
$ws.Range("M5").Value = 0.9494075
$ws.Range("N5").Value = 1.898815
$ws.Range("O5").Value = 0.002513328429799062
$ws.Range("P5").Value = 0.00167695720126921
$ws.Range("Q5").Value = 0.523155812355
$ws.Range("R5").Value = 3.13893487413
$ws.Range("S5").Value = 0.002513328429799062
$ws.Range("T5").Value = 0.00167695720126921

# Row 6
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.551034
$ws.Range("H6").Value = 1.653102
$ws.Range("M6").Value = 58.041021
$ws.Range("N6").Value = 174.123063
$ws.Range("O6").Value = 0.1536496690555577
$ws.Range("P6").Value = 0.1537785010150554
$ws.Range("Q6").Value = 31.982575965714
$ws.Range("R6").Value = 287.843183691426
$ws.Range("S6").Value = 0.1536496690555577
$ws.Range("T6").Value = 0.1537785010150554
